# ADD: prius CNG: citycar
#
# This script reproduces, on the data/content level, the changes shown by the
# OOXML diff for xl/worksheets/sheet1.xml (+ the shared strings it implies):
#
#   1) Rows 15-22 in column A are "uncommented": the leading "#" is removed
#      from "#citycar/dat/Japansese_city_car.dat" (rows 23-24 stay commented).
#   2) Rows 15-24 in column L: the png-folder reference is renamed from
#      "../images/OBKcars02" to "../images/OBKcars002".
#   3) C19 gets the Japanese-font style that C21:C24 already use.
#   4) Two brand new data rows (26 & 27) are appended describing a new
#      "prius CNG" city car ("TODOYA-Brius-1"/"TODOYA-Brius-2"), copied
#      in shape from the very first data row (row 5) and pointed at a new
#      image folder "../images/brius".
#   5) The active selection moves to C19 (and the view scrolls down a bit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Un-comment column A for rows 15-22 (rows 23 & 24 remain commented).
# ---------------------------------------------------------------------
for ($r = 15; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = "citycar/dat/Japansese_city_car.dat"
}

# ---------------------------------------------------------------------
# 2) Point column L (pngfile) at the renamed image folder for rows 15-24.
# ---------------------------------------------------------------------
for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 12).Value = "../images/OBKcars002"
}

# ---------------------------------------------------------------------
# 3) C19 picks up the Japanese-font style already used by C21:C24.
# ---------------------------------------------------------------------
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 4) Append the two new "prius CNG" city car rows (26 & 27), mirroring the
#    layout/style of the existing data rows.
# ---------------------------------------------------------------------
$ws.Range("A5:T5").Copy() | Out-Null
$ws.Range("A26:T26").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:T5").Copy() | Out-Null
$ws.Range("A27:T27").PasteSpecial(-4122) | Out-Null

# Row 26: TODOYA-Brius-1
$ws.Cells.Item(26, 1).Value = "citycar/dat/Japansese_city_car.dat"
$ws.Cells.Item(26, 2).ClearContents() | Out-Null
$ws.Cells.Item(26, 3).Value = "TODOYA-Brius-1"
$ws.Cells.Item(26, 4).Value = "citycar"
$ws.Cells.Item(26, 5).Value = "Harucarro"
$ws.Cells.Item(26, 6).Value = 135
$ws.Cells.Item(26, 7).Value = 11
$ws.Cells.Item(26, 8).Value = 2011
$ws.Cells.Item(26, 9).Value = 10
$ws.Cells.Item(26, 10).ClearContents() | Out-Null
$ws.Cells.Item(26, 11).ClearContents() | Out-Null
$ws.Cells.Item(26, 12).Value = "../images/brius"
$ws.Cells.Item(26, 13).Value = "'0.1"
$ws.Cells.Item(26, 14).Value = "'0.0"
$ws.Cells.Item(26, 15).Value = "'0.2"
$ws.Cells.Item(26, 16).Value = "'0.3"
$ws.Cells.Item(26, 17).Value = "'0.6"
$ws.Cells.Item(26, 18).Value = "'0.5"
$ws.Cells.Item(26, 19).Value = "'0.4"
$ws.Cells.Item(26, 20).Value = "'0.7"

# Row 27: TODOYA-Brius-2
$ws.Cells.Item(27, 1).Value = "citycar/dat/Japansese_city_car.dat"
$ws.Cells.Item(27, 2).ClearContents() | Out-Null
$ws.Cells.Item(27, 3).Value = "TODOYA-Brius-2"
$ws.Cells.Item(27, 4).Value = "citycar"
$ws.Cells.Item(27, 5).Value = "Harucarro"
$ws.Cells.Item(27, 6).Value = 130
$ws.Cells.Item(27, 7).Value = 9
$ws.Cells.Item(27, 8).Value = 2011
$ws.Cells.Item(27, 9).Value = 11
$ws.Cells.Item(27, 10).ClearContents() | Out-Null
$ws.Cells.Item(27, 11).ClearContents() | Out-Null
$ws.Cells.Item(27, 12).Value = "../images/brius"
$ws.Cells.Item(27, 13).Value = "'1.1"
$ws.Cells.Item(27, 14).Value = "'1.0"
$ws.Cells.Item(27, 15).Value = "'1.2"
$ws.Cells.Item(27, 16).Value = "'1.3"
$ws.Cells.Item(27, 17).Value = "'1.6"
$ws.Cells.Item(27, 18).Value = "'1.5"
$ws.Cells.Item(27, 19).Value = "'1.4"
$ws.Cells.Item(27, 20).Value = "'1.7"

# ---------------------------------------------------------------------
# 5) Move the selection/scroll position to C19.
# ---------------------------------------------------------------------
$ws.Range("C19").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 4
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
